# Update sequential diagram for "playlist del" command
# (mirrors commit "update sequential diagram for playlist del command")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) "playlist del p/Fav" -> "playlist del 1" (single-run textbox) ---
$sh = $s.Shapes.Item(12)
$sh.TextFrame.TextRange.Text = "playlist del 1"

# --- 2) execute("playlist del p/Fav") -> execute("playlist del 1") ---
$sh = $s.Shapes.Item(14)
$sh.TextFrame.TextRange.Text = "execute(“playlist del 1”)"

# --- 3) deletePlaylist(Fav) -> deletePlaylist(1)                      ---
#     Only the second run "(Fav)" changes; use Characters() so the
#     first run ("deletePlaylist") keeps its own run/formatting.
$sh = $s.Shapes.Item(16)
$tr = $sh.TextFrame.TextRange
$full = $tr.Text
$openParen = $full.IndexOf("(") + 1
$len = $full.Length - $openParen + 1
$tr.Characters($openParen, $len).Text = "(1)"

# --- 4) Date placeholder "10/31/18" -> "11/10/18" on the slide master
#        and every slide layout (both show the recalculated date stamp).
$m = $p.Slides.Item(1).Master
$m.Shapes.Item(3).TextFrame.TextRange.Text = "11/10/18"

for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
    $cl = $m.CustomLayouts.Item($i)
    $cl.Shapes.Item(1).TextFrame.TextRange.Text = "11/10/18"
}
